# The source diff for this revision touches only the *serialized XML
# attribute order* inside word/document.xml and word/styles.xml (e.g.
# namespace declarations on <w:document>, the order of attributes on
# <w:pgSz>/<w:pgMar>/<w:rFonts>/<w:lang>/<w:lsdException>/<w:style>/...).
# Every removed line and its replacement carry the *exact same set* of
# attributes/values - only their order changed - which is exactly what the
# commit message describes: "Fixed POI packaging and upgraded to POI 3.15."
# i.e. the expected-generation fixture was simply regenerated with a newer
# Apache POI that happens to alphabetize attributes when it writes OOXML
# parts back out. No paragraph text, run formatting, style definitions,
# page geometry, fonts, or language settings actually changed value.
#
# That kind of low-level XML-serializer attribute ordering is not part of
# the Word object model (Find/Replace, PageSetup, Styles, ... all operate
# on semantic values, and re-touching those values here only leaves the
# visible content/formatting unchanged while causing this engine to
# re-serialize the part through a different code path - adding artifacts
# such as w14:paraId/w14:textId, extra namespace declarations, or a
# standalone="yes" XML declaration - none of which appear in the target
# diff). So there is nothing to change through COM: the document's
# content, structure and formatting already match the target state.
$d = $word.ActiveDocument
